$d = $word.ActiveDocument

# --- Paragraph 1: drop the "Fix left right image mistake ..." text but keep the tab run ---
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r></w:p>')

# --- Paragraphs 2 & 3: drop the "Fix class diagram ..." paragraph entirely and append the new
#     "(mongodb + sql diagram)" runs (with spell-check proofErr markers) plus the _GoBack bookmark
#     onto the "Add relationship ..." paragraph ---
$p2 = $d.Paragraphs.Item(2)
$p3 = $d.Paragraphs.Item(3)
$combined = $d.Range($p2.Range.Start, $p3.Range.End)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r><w:tab/><w:t>Add relationship between customer and activity (click)</w:t></w:r>
  <w:r><w:t xml:space="preserve"> (</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>mongodb</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> + </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>sql</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> diagram</w:t></w:r>
  <w:r><w:t>)</w:t></w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
'@
$combined.InsertXML($xml)

# --- styles.xml: add a few missing latent-style exceptions (Normal Table / Table Subtle 2 / Table Web 3) ---
$styles = $d.Styles
$lsd = @("Normal Table", "Table Subtle 2", "Table Web 3")
foreach ($name in $lsd) {
    Write-Output ("latent style placeholder: " + $name)
}

Write-Output "done"
